{"js": "// The diff does two things to the body of the document:\n//  1. Shortens the long run of trailing spaces on the 3rd paragraph\n//     (the one ending in \"...tensorflow\" + spaces) from 90 down to 7.\n//  2. Appends five new paragraphs at the end of the document:\n//       - an empty paragraph\n//       - \"python -m pip install --upgrade --trusted-host pypi.org\n//          --trusted-host files.pythonhosted.org ipykernel -vvv\"\n//       - two empty paragraphs\n//       - \"ipython and ipykernel also update and uinstakkl \" + 83 spaces\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// --- 1. Shrink the trailing whitespace run on the 3rd paragraph ---------\nconst thirdParagraph = paragraphs.items[2];\nconst hits = thirdParagraph.search(\"tensorflow\", { matchCase: true });\nhits.load(\"text\");\nawait context.sync();\n\nconst afterWord = hits.items[0].getRange(\"After\");\nconst paragraphEnd = thirdParagraph.getRange(\"End\");\nconst trailingSpaces = afterWord.expandTo(paragraphEnd);\ntrailingSpaces.insertText(\"       \", \"Replace\"); // 7 spaces\n\n// --- 2. Append the new paragraphs ---------------------------------------\nlet cursor = body.paragraphs.getLast();\ncursor = cursor.insertParagraph(\"\", \"After\");\ncursor = cursor.insertParagraph(\n  \"python -m pip install --upgrade --trusted-host pypi.org --trusted-host files.pythonhosted.org ipykernel -vvv\",\n  \"After\"\n);\ncursor = cursor.insertParagraph(\"\", \"After\");\ncursor = cursor.insertParagraph(\"\", \"After\");\ncursor = cursor.insertParagraph(\n  \"ipython and ipykernel also update and uinstakkl \" + \" \".repeat(83),\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "# The diff does two things to the body of the document:\n#  1. Shortens the long run of trailing spaces on the 3rd paragraph\n#     (the one ending in \"...tensorflow\" + spaces) from 90 down to 7.\n#  2. Appends five new paragraphs at the end of the document:\n#       - an empty paragraph\n#       - \"python -m pip install --upgrade --trusted-host pypi.org\n#          --trusted-host files.pythonhosted.org ipykernel -vvv\"\n#       - two empty paragraphs\n#       - \"ipython and ipykernel also update and uinstakkl \" + 83 spaces\n\n$d = $word.ActiveDocument\n\n# --- 1. Shrink the trailing whitespace run on the 3rd paragraph -----------\n$p3 = $d.Paragraphs(3).Range\n$p3EndExclMark = $p3.End - 1            # exclude the paragraph mark itself\n$find = $p3.Find\n$find.Text = \"tensorflow\"\n$found = $find.Execute()                # $p3 collapses onto the match\n$trailing = $d.Range($p3.End, $p3EndExclMark)\n$trailing.Text = \"       \"              # 7 spaces\n\n# --- 2. Append the new paragraphs ------------------------------------------\nfunction NewPara([string]$text) {\n    $endRng = $d.Content\n    $endRng.Collapse(0)                 # wdCollapseEnd\n    $endRng.InsertParagraphAfter()\n    $lastPara = $d.Paragraphs($d.Paragraphs.Count).Range\n    $lastPara.Collapse(0)\n    if ($text.Length -gt 0) {\n        $lastPara.InsertBefore($text)\n    }\n}\n\n$spaces83 = \" \".PadLeft(83)\n\nNewPara \"\"\nNewPara \"python -m pip install --upgrade --trusted-host pypi.org --trusted-host files.pythonhosted.org ipykernel -vvv\"\nNewPara \"\"\nNewPara \"\"\nNewPara (\"ipython and ipykernel also update and uinstakkl \" + $spaces83)\n"}
